# [Fonds de solidarite] Add 2020-08-11 data
# Update nombre_aides (C) and montant_total (D) figures for rows reflecting
# the refreshed 2020-08-11 extract. Values are written as text (NumberFormat
# "@") to preserve the source workbook's inline-string cell typing, e.g. the
# trailing ".00" on whole-euro amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: nombre_aides 156 -> 159, montant_total 343600.00 -> 349600.00
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "159"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "349600.00"

# Row 3: nombre_aides 854 -> 869, montant_total 2327070.65 -> 2411251.17
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "869"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2411251.17"

# Row 4: nombre_aides 350 -> 355, montant_total 1236968.92 -> 1286368.92
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "355"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1286368.92"

# Row 5: nombre_aides 91 -> 93, montant_total 367982.09 -> 387982.09
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "93"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "387982.09"

# Row 33: nombre_aides 93 -> 96, montant_total 243826.00 -> 260748.00
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "96"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "260748.00"

# Row 34: nombre_aides 508 -> 519, montant_total 1563122.82 -> 1613326.26
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "519"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1613326.26"

# Row 35: nombre_aides 205 -> 210, montant_total 990347.11 -> 1027347.11
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "210"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1027347.11"

# Row 37: nombre_aides 23 -> 24, montant_total 140500.00 -> 150500.00
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "24"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "150500.00"

# Row 38: nombre_aides 19 -> 21, montant_total 42200.00 -> 46200.00
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "21"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "46200.00"

# Row 50: nombre_aides 93 -> 94, montant_total 257768.17 -> 267768.17
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "94"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "267768.17"

# Row 51: nombre_aides 541 -> 547, montant_total 1774758.52 -> 1830509.52
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "547"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1830509.52"

# Row 56: nombre_aides 668 -> 676, montant_total 1670796.41 -> 1690796.41
$ws.Range("C56").NumberFormat = "@"
$ws.Range("C56").Value = "676"
$ws.Range("D56").NumberFormat = "@"
$ws.Range("D56").Value = "1690796.41"

# Row 57: nombre_aides 3285 -> 3309, montant_total 9856825.10 -> 10010898.69
$ws.Range("C57").NumberFormat = "@"
$ws.Range("C57").Value = "3309"
$ws.Range("D57").NumberFormat = "@"
$ws.Range("D57").Value = "10010898.69"

# Row 58: nombre_aides 1695 -> 1702, montant_total 6745611.92 -> 6796111.92
$ws.Range("C58").NumberFormat = "@"
$ws.Range("C58").Value = "1702"
$ws.Range("D58").NumberFormat = "@"
$ws.Range("D58").Value = "6796111.92"

# Row 59: nombre_aides 578 -> 580, montant_total 2717640.96 -> 2737640.96
$ws.Range("C59").NumberFormat = "@"
$ws.Range("C59").Value = "580"
$ws.Range("D59").NumberFormat = "@"
$ws.Range("D59").Value = "2737640.96"

# Row 62: nombre_aides 272 -> 273, montant_total 638263.00 -> 640263.00
$ws.Range("C62").NumberFormat = "@"
$ws.Range("C62").Value = "273"
$ws.Range("D62").NumberFormat = "@"
$ws.Range("D62").Value = "640263.00"

# Row 72: nombre_aides 8 -> 12, montant_total 39000.00 -> 58000.00
$ws.Range("C72").NumberFormat = "@"
$ws.Range("C72").Value = "12"
$ws.Range("D72").NumberFormat = "@"
$ws.Range("D72").Value = "58000.00"

# Row 73: nombre_aides 84 -> 92, montant_total 218542.41 -> 242519.87
$ws.Range("C73").NumberFormat = "@"
$ws.Range("C73").Value = "92"
$ws.Range("D73").NumberFormat = "@"
$ws.Range("D73").Value = "242519.87"

# Row 74: nombre_aides 362 -> 377, montant_total 1067946.54 -> 1154173.74
$ws.Range("C74").NumberFormat = "@"
$ws.Range("C74").Value = "377"
$ws.Range("D74").NumberFormat = "@"
$ws.Range("D74").Value = "1154173.74"

# Row 75: nombre_aides 141 -> 147, montant_total 544092.18 -> 579358.18
$ws.Range("C75").NumberFormat = "@"
$ws.Range("C75").Value = "147"
$ws.Range("D75").NumberFormat = "@"
$ws.Range("D75").Value = "579358.18"

# Row 76: nombre_aides 39 -> 40, montant_total 167497.67 -> 177497.67
$ws.Range("C76").NumberFormat = "@"
$ws.Range("C76").Value = "40"
$ws.Range("D76").NumberFormat = "@"
$ws.Range("D76").Value = "177497.67"

# Row 77: nombre_aides 9 -> 10, montant_total 51000.00 -> 61000.00
$ws.Range("C77").NumberFormat = "@"
$ws.Range("C77").Value = "10"
$ws.Range("D77").NumberFormat = "@"
$ws.Range("D77").Value = "61000.00"

# Row 78: nombre_aides 8 -> 13, montant_total 16000.00 -> 26000.00
$ws.Range("C78").NumberFormat = "@"
$ws.Range("C78").Value = "13"
$ws.Range("D78").NumberFormat = "@"
$ws.Range("D78").Value = "26000.00"

# Row 80: nombre_aides 849 -> 851, montant_total 2613791.11 -> 2620997.11
$ws.Range("C80").NumberFormat = "@"
$ws.Range("C80").Value = "851"
$ws.Range("D80").NumberFormat = "@"
$ws.Range("D80").Value = "2620997.11"

# Row 81: nombre_aides 316 -> 318, montant_total 1220440.79 -> 1227440.79
$ws.Range("C81").NumberFormat = "@"
$ws.Range("C81").Value = "318"
$ws.Range("D81").NumberFormat = "@"
$ws.Range("D81").Value = "1227440.79"

# Row 82: nombre_aides 106 -> 107, montant_total 492484.52 -> 502484.52
$ws.Range("C82").NumberFormat = "@"
$ws.Range("C82").Value = "107"
$ws.Range("D82").NumberFormat = "@"
$ws.Range("D82").Value = "502484.52"

